# CRM: Auto-backup 2025-10-13 17:40:43
# Sets the "segundo_estatus" (column J) value for row 2 (customer C1000),
# which was previously blank, to "PEND.DOC.PARA EVALUACION".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Clientes")

$ws.Range("J2").Value = "PEND.DOC.PARA EVALUACION"
